# Updates crypto price/volume figures on Sheet1 (and restores the
# PancakeSwap / InternetComputer(DFINITY) row ordering for rows 30-31)
# to match the latest scrape, per the GitHub Actions cron job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.946.36'
$ws.Range("E2").Value = '  -1.36%  '
$ws.Range("D3").Value = '3.406.41'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.99'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.413'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '3.989.40'
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.98%  '
$ws.Range("D15").Value = '3.403.33'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '61.977.41'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.564'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '3.580.28'
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("E26").Value = '  -2.99%  '
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("E35").Value = '  +3.35%  '
$ws.Range("E36").Value = '  +2.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '31.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("D40").Value = '3.438.50'
$ws.Range("E40").Value = '  -1.79%  '
$ws.Range("E41").Value = '  +3.46%  '
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = '  -1.95%  '
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").Value = '2.544.02'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.05%  '
$ws.Range("E49").Value = '  -3.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("E51").Value = '  -0.02%  '
